# The deck's slide master currently carries the "Integral" theme
# (ppt/theme/theme1.xml), while its color values are being replaced with the
# stock "Office Theme" palette (the same palette that ppt/theme/theme2.xml,
# used by the notes master, already carries).
#
# PowerPoint doesn't expose a single "swap these two theme parts" verb, so we
# drive this the way a user actually would from the Slide Master view: select
# the master and recolor its theme, one theme color at a time, via
# ThemeColorScheme.
#
# ThemeColorScheme.Colors(n) is indexed using the standard
# msoThemeColorSchemeIndex order:
#   1 Dark1, 2 Light1, 3 Dark2, 4 Light2,
#   5 Accent1 .. 10 Accent6, 11 Hyperlink, 12 FollowedHyperlink
# and each ColorFormat's .RGB takes a packed BGR long (R + G*256 + B*65536),
# exactly like the classic VBA RGB() function.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$themeColors = $theme.ThemeColorScheme

# Target "Office Theme" color scheme, expressed as packed RGB longs.
$officeThemeRgb = @{
    1  = 0          # Dark1             000000
    2  = 16777215   # Light1            FFFFFF
    3  = 6968388     # Dark2             44546A
    4  = 15132391    # Light2            E7E6E6
    5  = 13998939    # Accent1           5B9BD5
    6  = 3243501     # Accent2           ED7D31
    7  = 10855845    # Accent3           A5A5A5
    8  = 49407       # Accent4           FFC000
    9  = 12874308    # Accent5           4472C4
    10 = 4697456      # Accent6           70AD47
    11 = 12673797     # Hyperlink         0563C1
    12 = 7491477      # FollowedHyperlink 954F72
}

for ($i = 1; $i -le 12; $i++) {
    $themeColors.Colors($i).RGB = $officeThemeRgb[$i]
}

# Best-effort rename to match the target theme/color-scheme naming.
$theme.Name = "Office Theme"
$themeColors.Name = "Office"
